$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cells that changed from "W" to "WS"
$ws.Range("G3").Value = "WS"
$ws.Range("O3").Value = "WS"
$ws.Range("C7").Value = "WS"
$ws.Range("T15").Value = "WS"
$ws.Range("E20").Value = "WS"
$ws.Range("P20").Value = "WS"

# Update the active selection to C7
$ws.Range("C7").Select()
